# 06.DanhSachChucNang.xlsx
#  - Chinh sua store select nhan vien: di chuyen vung nhin / vung chon hien
#    tai cua Sheet1 (view scrolled to A31 / selection E44 -> scrolled to A4 /
#    selection E15).
#  - Cap nhat danh sach chuc nang: dien % hoan thanh con thieu cho cac dong
#    9 (E12), 21 (E24) va 22 (E25).
#  - Them validator gia ve cho wucThemTuyen: du lieu "Phan cong" cua cac dong
#    21/22 bi go nham vao cot G (ngoai bang Table1, A3:F48) - dua ve dung cot
#    F "Phan cong" trong bang.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Cap nhat danh sach chuc nang: % Hoan thanh con trong -----------------

# Dong 9 (STT=9, Usecase 2.4 - "Xem phan hoi khach hang")
$ws.Range("E12").Value = 0.9

# Dong 21 (STT=21, Usecase 4.7)
$ws.Range("E24").Value = 0.9

# Dong 22 (STT=22, Usecase 4.8)
$ws.Range("E25").Value = 1

# --- Them validator gia ve cho wucThemTuyen: sua cot "Phan cong" ----------
# Gia tri "Tu" bi nhap lech sang cot G (ngoai Table1) o 2 dong 21 va 22;
# chuyen lai vao cot F ("Phan cong") va xoa o G.

$ws.Range("F24").Value = $ws.Range("G24").Value2
$ws.Range("G24").Clear()

$ws.Range("F25").Value = $ws.Range("G25").Value2
$ws.Range("G25").Clear()

# --- Chinh sua store select nhan vien: cap nhat vung nhin / vung chon -----

$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("E15").Select()
